$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns, matching the existing header style (bold, bordered, centered)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins/Losses/Ties) for every player row
$ws.Range("AD2:AD71").Value = 71
$ws.Range("AE2:AE71").Value = 91
$ws.Range("AF2:AF71").Value = 0
